$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# --- Rows 2-36: update Price (D) / Volume(1h) (E) only ---
Set-TextValue $ws.Range('D2') '30.775.39'
Set-TextValue $ws.Range('E2') '  -1.34%  '
Set-TextValue $ws.Range('D3') '1.940.78'
Set-TextValue $ws.Range('E3') '  -0.99%  '
Set-TextValue $ws.Range('E4') '  +0.19%  '
Set-TextValue $ws.Range('D5') '242.05'
Set-TextValue $ws.Range('E5') '  -2.06%  '
Set-TextValue $ws.Range('E6') '  +0.20%  '
Set-TextValue $ws.Range('D7') '0.4884'
Set-TextValue $ws.Range('E7') '  -0.18%  '
Set-TextValue $ws.Range('D8') '0.2935'
Set-TextValue $ws.Range('E8') '  -1.22%  '
Set-TextValue $ws.Range('D9') '0.06928'
Set-TextValue $ws.Range('E9') '  +1.29%  '
Set-TextValue $ws.Range('D10') '19.47'
Set-TextValue $ws.Range('E10') '  +1.54%  '
Set-TextValue $ws.Range('D11') '105.97'
Set-TextValue $ws.Range('E11') '  -0.39%  '
Set-TextValue $ws.Range('D12') '1.938.82'
Set-TextValue $ws.Range('E12') '  -0.85%  '
Set-TextValue $ws.Range('D13') '0.07745'
Set-TextValue $ws.Range('E13') '  -0.32%  '
Set-TextValue $ws.Range('D14') '5.350'
Set-TextValue $ws.Range('E14') '  -1.36%  '
Set-TextValue $ws.Range('E15') '  -1.74%  '
Set-TextValue $ws.Range('D16') '275.40'
Set-TextValue $ws.Range('E16') '  -4.16%  '
Set-TextValue $ws.Range('D17') '30.782.14'
Set-TextValue $ws.Range('E17') '  -1.34%  '
Set-TextValue $ws.Range('E18') '  -0.33%  '
Set-TextValue $ws.Range('D19') '13.13'
Set-TextValue $ws.Range('E19') '  -0.89%  '
Set-TextValue $ws.Range('D20') '2.203.08'
Set-TextValue $ws.Range('E20') '  +0.38%  '
Set-TextValue $ws.Range('D21') '1.000'
Set-TextValue $ws.Range('E21') '  +0.08%  '
Set-TextValue $ws.Range('D22') '5.435'
Set-TextValue $ws.Range('E22') '  -2.78%  '
Set-TextValue $ws.Range('D23') '1.002'
Set-TextValue $ws.Range('E23') '  +0.85%  '
Set-TextValue $ws.Range('D24') '6.503'
Set-TextValue $ws.Range('E24') '  -1.35%  '
Set-TextValue $ws.Range('D25') '9.725'
Set-TextValue $ws.Range('E25') '  -2.71%  '
Set-TextValue $ws.Range('D26') '169.09'
Set-TextValue $ws.Range('E26') '  +0.49%  '
Set-TextValue $ws.Range('D27') '19.60'
Set-TextValue $ws.Range('E27') '  -2.09%  '
Set-TextValue $ws.Range('D28') '2.159'
Set-TextValue $ws.Range('E28') '  -2.07%  '
Set-TextValue $ws.Range('D29') '0.1040'
Set-TextValue $ws.Range('E29') '  -2.43%  '
Set-TextValue $ws.Range('D30') '1.388'
Set-TextValue $ws.Range('E30') '  -3.72%  '
Set-TextValue $ws.Range('D32') '4.557'
Set-TextValue $ws.Range('E32') '  -4.73%  '
Set-TextValue $ws.Range('E33') '  -2.94%  '
Set-TextValue $ws.Range('D34') '0.04861'
Set-TextValue $ws.Range('E34') '  -3.51%  '
Set-TextValue $ws.Range('D35') '0.7499'
Set-TextValue $ws.Range('E35') '  -2.46%  '
Set-TextValue $ws.Range('D36') '1.156'
Set-TextValue $ws.Range('E36') '  -0.86%  '

# --- Rows 37-51: new "Frax" row inserted, pushing everything down one; Elrond drops off ---
Set-TextValue $ws.Range('B37') 'Frax'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D37') '1.000'
Set-TextValue $ws.Range('E37') '  +0.22%  '
Set-TextValue $ws.Range('B38') 'HuobiToken'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D38') '2.726'
Set-TextValue $ws.Range('E38') '  -0.04%  '
Set-TextValue $ws.Range('B39') 'VeChain'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D39') '0.01990'
Set-TextValue $ws.Range('E39') '  -3.01%  '
Set-TextValue $ws.Range('B40') 'MXToken'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D40') '2.665'
Set-TextValue $ws.Range('E40') '  -1.80%  '
Set-TextValue $ws.Range('B41') 'Aave'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D41') '77.92'
Set-TextValue $ws.Range('E41') '  +6.05%  '
Set-TextValue $ws.Range('B42') 'FraxShare'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D42') '6.466'
Set-TextValue $ws.Range('E42') '  +0.71%  '
Set-TextValue $ws.Range('B43') 'RenderToken'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D43') '2.105'
Set-TextValue $ws.Range('E43') '  -1.20%  '
Set-TextValue $ws.Range('B44') 'TrustWalletToken'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D44') '0.9011'
Set-TextValue $ws.Range('E44') '  +1.68%  '
Set-TextValue $ws.Range('B45') 'Quant'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D45') '108.66'
Set-TextValue $ws.Range('E45') '  -1.05%  '
Set-TextValue $ws.Range('B46') 'TheSandbox'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D46') '0.4405'
Set-TextValue $ws.Range('E46') '  -1.19%  '
Set-TextValue $ws.Range('B47') 'PaxDollar'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D47') '0.9973'
Set-TextValue $ws.Range('E47') '  -0.13%  '
Set-TextValue $ws.Range('B48') 'Aptos'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D48') '7.727'
Set-TextValue $ws.Range('E48') '  +3.19%  '
Set-TextValue $ws.Range('B49') 'Maker'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D49') '990.62'
Set-TextValue $ws.Range('E49') '  -0.15%  '
Set-TextValue $ws.Range('B50') 'Algorand'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D50') '0.1245'
Set-TextValue $ws.Range('E50') '  -2.02%  '
Set-TextValue $ws.Range('B51') 'EnergySwap'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D51') '9.236'
Set-TextValue $ws.Range('E51') '  -1.81%  '
